$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the test import values in row 2 (G2:I2) which previously held
# placeholder/English strings and now hold the real French labels.
$ws.Range("G2").Value = "Restaurant avec cuisine sur place"
$ws.Range("H2").Value = "Concédée"
$ws.Range("I2").Value = "Public"

# G2 keeps its text format but switches to the Times New Roman font used
# elsewhere in the sheet.
$ws.Range("G2").Font.Name = "Times New Roman"

# H2 and I2 switch from the forced text format to General and wrap their
# text so the longer French labels display properly.
$ws.Range("H2:I2").NumberFormat = "General"
$ws.Range("H2:I2").WrapText = $true

# The active selection moves from F2 to I2.
[void]$ws.Range("I2").Select()
